$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''69.868.29'
$ws.Range("E2").Value = '  -0.85%  '
$ws.Range("D3").Value = '''3.538.72'
$ws.Range("E3").Value = '  -1.10%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '''611.58'
$ws.Range("E5").Value = '  +4.33%  '
$ws.Range("D6").Value = '''184.89'
$ws.Range("E6").Value = '  -0.57%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").Value = '''0.215'
$ws.Range("E9").Value = '  +2.62%  '
$ws.Range("D10").Value = '''0.644'
$ws.Range("E10").Value = '  -0.82%  '
$ws.Range("E11").Value = '  -1.33%  '
$ws.Range("E12").Value = '  -2.44%  '
$ws.Range("D13").Value = '''9.42'
$ws.Range("E13").Value = '  -1.32%  '
$ws.Range("D14").Value = '''4.101.42'
$ws.Range("E14").Value = '  -1.07%  '
$ws.Range("D15").Value = '''609.73'
$ws.Range("E15").Value = '  +7.68%  '
$ws.Range("D16").Value = '''69.916.25'
$ws.Range("E16").Value = '  -0.74%  '
$ws.Range("D17").Value = '''3.549.24'
$ws.Range("E17").Value = '  -0.27%  '
$ws.Range("D18").Value = '''12.64'
$ws.Range("E18").Value = '  +1.89%  '
$ws.Range("D19").Value = '''18.82'
$ws.Range("E19").Value = '  -3.92%  '
$ws.Range("E20").Value = '  -0.28%  '
$ws.Range("E21").Value = '  -2.38%  '
$ws.Range("E22").Value = '  -1.52%  '
$ws.Range("E23").Value = '  +1.59%  '
$ws.Range("D24").Value = '''100.17'
$ws.Range("E24").Value = '  +5.06%  '
$ws.Range("E25").Value = '  +0.33%  '
$ws.Range("E26").Value = '  +1.07%  '
$ws.Range("E27").Value = '  -5.31%  '
$ws.Range("D28").Value = '''9.58'
$ws.Range("E28").Value = '  +4.80%  '
$ws.Range("D29").Value = '''32.43'
$ws.Range("E29").Value = '  +0.63%  '
$ws.Range("D30").Value = '''7.01'
$ws.Range("E30").Value = '  -4.31%  '
$ws.Range("D31").Value = '''12.22'
$ws.Range("E31").Value = '  -2.04%  '
$ws.Range("E32").Value = '  -1.05%  '
$ws.Range("D33").Value = '''63.41'
$ws.Range("E33").Value = '  -2.23%  '
$ws.Range("D34").Value = '''3.60'
$ws.Range("E34").Value = '  +18.36%  '
$ws.Range("D35").Value = '''3.24'
$ws.Range("E35").Value = '  -3.39%  '
$ws.Range("D36").Value = '''533.88'
$ws.Range("E36").Value = '  -5.58%  '
$ws.Range("B37").Value = 'Dai'
$ws.Range("C37").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D37").Value = '''1.00'
$ws.Range("E37").Value = '  -0.04%  '
$ws.Range("B38").Value = 'TheGraph'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D38").Value = '''0.400'
$ws.Range("E38").Value = '  -5.17%  '
$ws.Range("D39").Value = '''37.04'
$ws.Range("E39").Value = '  -2.10%  '
$ws.Range("D40").Value = '''0.0₃0780'
$ws.Range("E40").Value = '  +0.43%  '
$ws.Range("D41").Value = '''3.56'
$ws.Range("E41").Value = '  +5.47%  '
$ws.Range("D42").Value = '''3.533.75'
$ws.Range("E42").Value = '  +4.89%  '
$ws.Range("D43").Value = '''0.137'
$ws.Range("E43").Value = '  +1.49%  '
$ws.Range("E44").Value = '  +2.52%  '
$ws.Range("E45").Value = '  -1.63%  '
$ws.Range("E46").Value = '  +3.97%  '
$ws.Range("E47").Value = '  -5.08%  '
$ws.Range("D48").Value = '''9.15'
$ws.Range("E48").Value = '  -3.12%  '
$ws.Range("E49").Value = '  +0.31%  '
$ws.Range("E50").Value = '  -3.32%  '
$ws.Range("D51").Value = '''135.73'
$ws.Range("E51").Value = '  -1.32%  '

# Clear the auto-applied quote-prefix formatting so styling matches the source
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D51").Style = "Normal"
